$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 64, pushing the existing data (old rows 64-114) down to 65-115.
$ws.Rows.Item(64).Insert()

# Fill in the new weekly price-report record for row 64.
$ws.Cells.Item(64, 1).Value = 6
$ws.Cells.Item(64, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(64, 3).Value = "Metropolitana"
$ws.Cells.Item(64, 4).Value = 44977
$ws.Cells.Item(64, 5).Value = 13
$ws.Cells.Item(64, 6).Value = 100114007
$ws.Cells.Item(64, 7).Value = "Jengibre"
$ws.Cells.Item(64, 8).Value = "Sin especificar"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 400
$ws.Cells.Item(64, 11).Value = 12000
$ws.Cells.Item(64, 12).Value = 13000
$ws.Cells.Item(64, 13).Value = 12425
$ws.Cells.Item(64, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(64, 15).Value = "Perú"
$ws.Cells.Item(64, 16).Value = 956
$ws.Cells.Item(64, 17).Value = 13
$ws.Cells.Item(64, 18).Value = "Hortaliza"
